$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BI")

# --- H2 trade-link matrix header (row 12) ---
# Previously C12 was blank with a note in D12 saying the H2 matrix didn't
# work yet; now the matrix is filled in, so the row 12 header mirrors row 3
# ("~TradeLinks") and the note is removed.
$ws.Range("C3").Copy()
$ws.Range("C12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C12").Value = "~TradeLinks"
$ws.Range("D12").ClearContents()

# --- Fill in the H2 trade-link matrix values (rows 14-19) ---
# Each new "1" cell gets the same format as the other matrix flag cells
# (gray fill + border), then the value is set.
$ws.Range("E5").Copy()

$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = 1

$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1

$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 1

$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I15").Value = 1

$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 1

$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = 1

$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 1

$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = 1

# --- Selection moved to L12 (matches the author's final cursor position) ---
$ws.Range("L12").Select()
